$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1): rename columns and reorder metric headers ---
$ws.Range("A1").Value = "Paso"
$ws.Range("B1").Value = "Config"
$ws.Range("C1").Value = "Dist"
$ws.Range("D1").Value = "Var"
$ws.Range("E1").Value = "Diferenciacion"
$ws.Range("F1").Value = "Block Bootstrapping"
$ws.Range("G1").Value = "Sieve Bootstrap"
$ws.Range("H1").Value = "LSPM"
$ws.Range("I1").Value = "LSPMW"
$ws.Range("J1").Value = "AREPD"
$ws.Range("K1").Value = "MCPS"
$ws.Range("L1").Value = "AV-MCPS"
$ws.Range("M1").Value = "DeepAR"
$ws.Range("N1").Value = "EnCQR-LSTM"

# --- Update data rows 2-25: reorder/replace metric columns F..N ---
# Build a rectangular array (24 rows x 9 cols) for columns F:N, rows 2:25
$data = New-Object 'object[,]' 24,9
$data[0,0] = 0.7126554218325974; $data[0,1] = 0.5766501968783261; $data[0,2] = 0.5782300834920086; $data[0,3] = 0.5732512545132614; $data[0,4] = 0.7057938330819035; $data[0,5] = 0.6291257877327282; $data[0,6] = 0.6716091734854959; $data[0,7] = 0.5706096009454333; $data[0,8] = 1.112709936422421
$data[1,0] = 0.8743521934894305; $data[1,1] = 0.5853046480970127; $data[1,2] = 1.450983088048725; $data[1,3] = 1.324120047199519; $data[1,4] = 0.8908609375329961; $data[1,5] = 0.6948233016659645; $data[1,6] = 0.6705485377379639; $data[1,7] = 0.588298910317083; $data[1,8] = 1.059723850405045
$data[2,0] = 0.6324559014148136; $data[2,1] = 0.5659788120944055; $data[2,2] = 0.8031842229937988; $data[2,3] = 0.8784972426811248; $data[2,4] = 0.6264654340791795; $data[2,5] = 0.5964491206202314; $data[2,6] = 0.6043629777114435; $data[2,7] = 0.5675513096607806; $data[2,8] = 1.022809791875938
$data[3,0] = 0.5821791912275598; $data[3,1] = 0.5715087696534736; $data[3,2] = 0.5597737657001477; $data[3,3] = 0.7827440109487351; $data[3,4] = 0.5857009941903075; $data[3,5] = 0.583579573183647; $data[3,6] = 0.5908712941750363; $data[3,7] = 0.571989034963788; $data[3,8] = 1.026842104055089
$data[4,0] = 0.5892182285019585; $data[4,1] = 0.5755720121442439; $data[4,2] = 0.5746211706825725; $data[4,3] = 0.6419531616488163; $data[4,4] = 0.5774599606081097; $data[4,5] = 0.7460416780380763; $data[4,6] = 1.121806008118179; $data[4,7] = 0.5901454652413448; $data[4,8] = 1.051972420586647
$data[5,0] = 0.94935715741024; $data[5,1] = 0.6495094704620347; $data[5,2] = 1.037753543256499; $data[5,3] = 0.7068784087372422; $data[5,4] = 0.9907819881634593; $data[5,5] = 1.346794202830598; $data[5,6] = 1.278431621376492; $data[5,7] = 0.6157641502706026; $data[5,8] = 1.202781442064724
$data[6,0] = 0.6547457039131218; $data[6,1] = 0.5628138726880555; $data[6,2] = 1.527809472419873; $data[6,3] = 0.9105795976079453; $data[6,4] = 0.6365520151925801; $data[6,5] = 0.581980425990367; $data[6,6] = 0.5841508683063711; $data[6,7] = 0.6118959602258814; $data[6,8] = 1.023943091723387
$data[7,0] = 0.7200803105868399; $data[7,1] = 0.6106961088801163; $data[7,2] = 0.5558338967728091; $data[7,3] = 1.092977909049395; $data[7,4] = 0.7311125363119273; $data[7,5] = 0.6289371268745428; $data[7,6] = 0.6133816891042274; $data[7,7] = 0.5548263216522614; $data[7,8] = 1.03022262859112
$data[8,0] = 0.5991784682094614; $data[8,1] = 0.6698427334953047; $data[8,2] = 0.7236462043541498; $data[8,3] = 0.7313423899471424; $data[8,4] = 0.5903121829430349; $data[8,5] = 0.6701350004923454; $data[8,6] = 0.6802083088175948; $data[8,7] = 0.623348094815062; $data[8,8] = 1.040907702819046
$data[9,0] = 0.990401152658249; $data[9,1] = 0.6368517624697899; $data[9,2] = 1.094303761026129; $data[9,3] = 1.37980964515669; $data[9,4] = 0.9186681073080789; $data[9,5] = 0.7760188156263539; $data[9,6] = 0.7133174782481635; $data[9,7] = 0.5772614482446846; $data[9,8] = 1.063541678902269
$data[10,0] = 0.593806997479187; $data[10,1] = 0.6176357044147445; $data[10,2] = 1.209632797065104; $data[10,3] = 0.6145172404111509; $data[10,4] = 0.5999768077750933; $data[10,5] = 0.666094371012726; $data[10,6] = 0.6240949649974504; $data[10,7] = 0.6835724464738447; $data[10,8] = 1.064609034910676
$data[11,0] = 0.596173703582406; $data[11,1] = 0.693760662127908; $data[11,2] = 0.5494253916601892; $data[11,3] = 0.5608265051388281; $data[11,4] = 0.606197133342542; $data[11,5] = 0.6433063020586817; $data[11,6] = 0.563674167455626; $data[11,7] = 0.5572324955172964; $data[11,8] = 1.079160430095435
$data[12,0] = 0.5668329692763; $data[12,1] = 0.5672070965782691; $data[12,2] = 0.5699896497385639; $data[12,3] = 0.57791362570924; $data[12,4] = 0.607241273656135; $data[12,5] = 0.749507124632801; $data[12,6] = 0.63195843067441; $data[12,7] = 0.5649703935124497; $data[12,8] = 0.8173496813035664
$data[13,0] = 0.8011668737375255; $data[13,1] = 0.5840019402228158; $data[13,2] = 0.8202640749220081; $data[13,3] = 0.9037813139279145; $data[13,4] = 0.9599374597024464; $data[13,5] = 0.8110787018069477; $data[13,6] = 0.8691802955043432; $data[13,7] = 0.5762861103346959; $data[13,8] = 0.8989764503656629
$data[14,0] = 0.5570378508949394; $data[14,1] = 0.5562894444892267; $data[14,2] = 0.6630046900746841; $data[14,3] = 0.5702952467489729; $data[14,4] = 0.6032079110748835; $data[14,5] = 0.7929706989746772; $data[14,6] = 0.5615837950447297; $data[14,7] = 0.5679141985543038; $data[14,8] = 0.8006498799270003
$data[15,0] = 0.7250873468229858; $data[15,1] = 0.5789101827732127; $data[15,2] = 0.6621145821756242; $data[15,3] = 0.7727081414334852; $data[15,4] = 0.8127193066840899; $data[15,5] = 0.7135535098756202; $data[15,6] = 0.7449171295124811; $data[15,7] = 0.582111051193741; $data[15,8] = 0.8404046955762926
$data[16,0] = 0.5762970435747076; $data[16,1] = 0.557951144576361; $data[16,2] = 0.7296249429994273; $data[16,3] = 0.5605753575630967; $data[16,4] = 0.5957741585919769; $data[16,5] = 0.5991713305275943; $data[16,6] = 0.5717296204259081; $data[16,7] = 0.5610235942103873; $data[16,8] = 0.820027559234253
$data[17,0] = 0.9963736355967837; $data[17,1] = 0.588294285548589; $data[17,2] = 0.8812946174388349; $data[17,3] = 0.9555740919234618; $data[17,4] = 1.083988254050418; $data[17,5] = 0.8743090602556495; $data[17,6] = 1.10847823222863; $data[17,7] = 0.6274855290585313; $data[17,8] = 1.055027373882682
$data[18,0] = 0.5727403929211751; $data[18,1] = 0.5750724248912124; $data[18,2] = 1.143535915599847; $data[18,3] = 0.5995390386705647; $data[18,4] = 0.6244534860332279; $data[18,5] = 0.6393403100645219; $data[18,6] = 0.6495013231615474; $data[18,7] = 0.5933760959151286; $data[18,8] = 0.8086248997392134
$data[19,0] = 0.635254873444856; $data[19,1] = 0.5747654680268093; $data[19,2] = 0.6530831576429931; $data[19,3] = 0.6108145954973647; $data[19,4] = 0.6627571709779209; $data[19,5] = 0.6159788096637256; $data[19,6] = 0.686165127218256; $data[19,7] = 0.5936987325509109; $data[19,8] = 0.8576438121889344
$data[20,0] = 0.6166640517859907; $data[20,1] = 0.5639203973202785; $data[20,2] = 0.5588899200459683; $data[20,3] = 0.6013270098371213; $data[20,4] = 0.6497496606752713; $data[20,5] = 0.638059074498138; $data[20,6] = 0.6394395803142389; $data[20,7] = 0.5611619270496143; $data[20,8] = 0.8515393028359909
$data[21,0] = 0.7329513392486269; $data[21,1] = 0.575849999177136; $data[21,2] = 0.9023532226570946; $data[21,3] = 0.8043491720067202; $data[21,4] = 0.8480261606611101; $data[21,5] = 1.154884489308201; $data[21,6] = 0.9502746002678047; $data[21,7] = 0.5873257670379519; $data[21,8] = 0.8596102923638532
$data[22,0] = 0.6116684245182376; $data[22,1] = 0.5808991827055636; $data[22,2] = 0.8318908898490085; $data[22,3] = 0.5974550344492254; $data[22,4] = 0.6434434589628736; $data[22,5] = 0.6830234903337824; $data[22,6] = 0.7252431028525861; $data[22,7] = 0.5891459325951937; $data[22,8] = 0.8393407005799141
$data[23,0] = 0.5962876166093992; $data[23,1] = 0.5779684890189307; $data[23,2] = 0.7198643536281423; $data[23,3] = 0.6357867839450184; $data[23,4] = 0.6661136382177197; $data[23,5] = 0.6866654549183088; $data[23,6] = 1.097018854895182; $data[23,7] = 0.5952914921738442; $data[23,8] = 0.8218890827182036

$ws.Range("F2:N25").Value = $data

# --- Remove now-unused column O entirely (content + formatting) ---
$ws.Range("O1:O25").Clear()
